# Probe interface change: append 5 new alarm rows (rows 4-8) to Sheet1,
# matching the new probe payload shape coming from the device integration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A/B can contain values that look numeric/date-like ("46",
# "2024-10-19"). Force those cells to Text format *before* assigning the
# values so Excel keeps them as literal strings instead of silently
# coercing them into numbers / date serials.
$ws.Range("A4:B8").NumberFormat = "@"

# --- Row 4: alarm #2, Heart Rate / Low (no Date/Participant/Block yet) ---
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "Heart Rate"
$ws.Range("F4").Value = "Low"
$ws.Range("G4").Value = "2024-10-19T12:06:32.829"

# --- Row 5: alarm #2, Blood Pressure / High ---
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "Blood Pressure"
$ws.Range("F5").Value = "High"
$ws.Range("G5").Value = "2024-10-19T12:06:33.660"

# --- Row 6: alarm #2, Oxygen Saturation / Very Low ---
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "Oxygen Saturation"
$ws.Range("F6").Value = "Very Low"
$ws.Range("G6").Value = "2024-10-19T12:06:34.831"

# --- Row 7: alarm #3, Heart Rate / Very Low, new Date/Participant/Block ---
$ws.Range("A7").Value = "2024-10-19"
$ws.Range("B7").Value = "46"
$ws.Range("C7").Value = "test3"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = "Heart Rate"
$ws.Range("F7").Value = "Very Low"
$ws.Range("G7").Value = "2024-10-19T12:12:04.549"

# --- Row 8: alarm #3, Blood Pressure / High ---
$ws.Range("A8").Value = "2024-10-19"
$ws.Range("B8").Value = "46"
$ws.Range("C8").Value = "test3"
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = "Blood Pressure"
$ws.Range("F8").Value = "High"
$ws.Range("G8").Value = "2024-10-19T12:12:05.137"
